$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 347.25
$ws.Range("I55").Value = 347.25
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 347.25
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -133.25
$ws.Range("N55").ClearContents()
$ws.Range("H94").Value = 5000
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H107").Value = 597.4
$ws.Range("I107").Value = 605.4286
$ws.Range("J107").Value = 485
$ws.Range("K107").Value = 605.4286
$ws.Range("L107").Value = 485
$ws.Range("M107").Value = 1314.5714
$ws.Range("N107").Value = -4325
$ws.Range("H116").Value = 5104.7896
$ws.Range("I116").Value = 3536.875
$ws.Range("J116").Value = 6245.091
$ws.Range("K116").Value = 3536.875
$ws.Range("L116").Value = 6245.091
$ws.Range("M116").Value = -94.875
$ws.Range("N116").Value = -13129.091
$ws.Range("H129").Value = 1271.3182
$ws.Range("I129").Value = 482.57144
$ws.Range("J129").Value = 1639.4
$ws.Range("K129").Value = 1447.71432
$ws.Range("L129").Value = 4918.200000000001
$ws.Range("M129").Value = 3552.28568
$ws.Range("N129").Value = -14918.2
$ws.Range("H137").Value = 2633467.2
$ws.Range("I137").Value = 2150.4092
$ws.Range("J137").Value = 6251528
$ws.Range("K137").Value = 6451.2276
$ws.Range("L137").Value = 18754584
$ws.Range("M137").Value = -3901.2276
$ws.Range("N137").Value = -18759684
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2499.6956
$ws.Range("I45").Value = 1639.5385
$ws.Range("J45").Value = 3617.9
$ws.Range("K45").Value = 1639.5385
$ws.Range("L45").Value = 3617.9
$ws.Range("M45").Value = -1262.5385
$ws.Range("N45").Value = -4371.9
$ws.Range("H102").Value = 2050.1667
$ws.Range("I102").Value = 2050.1667
$ws.Range("K102").Value = 2050.1667
$ws.Range("M102").Value = -428.1667000000002
$ws.Range("H122").Value = 1976.6154
$ws.Range("I122").Value = 1887.75
$ws.Range("J122").Value = 2118.8
$ws.Range("K122").Value = 5663.25
$ws.Range("L122").Value = 6356.400000000001
$ws.Range("M122").Value = -3213.25
$ws.Range("N122").Value = -11256.4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2167.5881
$ws.Range("I20").Value = 1261.7142
$ws.Range("J20").Value = 2801.7
$ws.Range("K20").Value = 1261.7142
$ws.Range("L20").Value = 2801.7
$ws.Range("M20").Value = -1014.7142
$ws.Range("N20").Value = -3295.7
$ws.Range("H86").Value = 1746.5
$ws.Range("I86").Value = 1399.6666
$ws.Range("J86").Value = 2093.3333
$ws.Range("K86").Value = 1399.6666
$ws.Range("L86").Value = 2093.3333
$ws.Range("M86").Value = -276.6666
$ws.Range("N86").Value = -4339.3333
$ws.Range("H89").Value = 1746.5
$ws.Range("I89").Value = 1399.6666
$ws.Range("J89").Value = 2093.3333
$ws.Range("K89").Value = 6998.333000000001
$ws.Range("L89").Value = 10466.6665
$ws.Range("M89").Value = -1382.333000000001
$ws.Range("N89").Value = -21698.6665
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6174632.5
$ws.Range("I31").Value = 1339.122
$ws.Range("K31").Value = 1339.122
$ws.Range("M31").Value = -1044.122
$ws.Range("H34").Value = 6174632.5
$ws.Range("I34").Value = 1339.122
$ws.Range("K34").Value = 1339.122
$ws.Range("M34").Value = -1137.122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 3728.1875
$ws.Range("I22").Value = 5525.5
$ws.Range("J22").Value = 3471.4285
$ws.Range("K22").Value = 16576.5
$ws.Range("L22").Value = 10414.2855
$ws.Range("M22").Value = -16407.5
$ws.Range("N22").Value = -10752.2855
$ws.Range("H27").Value = 3728.1875
$ws.Range("I27").Value = 5525.5
$ws.Range("J27").Value = 3471.4285
$ws.Range("K27").Value = 16576.5
$ws.Range("L27").Value = 10414.2855
$ws.Range("M27").Value = -16474.5
$ws.Range("N27").Value = -10618.2855
$ws.Range("H81").Value = 2527.6667
$ws.Range("I81").Value = 833.3333
$ws.Range("J81").Value = 2951.25
$ws.Range("K81").Value = 2499.9999
$ws.Range("L81").Value = 8853.75
$ws.Range("M81").Value = -1376.9999
$ws.Range("N81").Value = -11099.75
$ws.Range("H84").Value = 2527.6667
$ws.Range("I84").Value = 833.3333
$ws.Range("J84").Value = 2951.25
$ws.Range("K84").Value = 7499.9997
$ws.Range("L84").Value = 26561.25
$ws.Range("M84").Value = -1883.9997
$ws.Range("N84").Value = -37793.25
$ws.Range("H131").Value = 263964.03
$ws.Range("J131").Value = 1121.762
$ws.Range("L131").Value = 3365.286
$ws.Range("N131").Value = -13445.286
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 15000
$ws.Range("J74").Value = 15000
$ws.Range("L74").Value = 15000
$ws.Range("N74").Value = -16996
$ws.Range("H77").Value = 15000
$ws.Range("J77").Value = 15000
$ws.Range("L77").Value = 45000
$ws.Range("N77").Value = -54984
$ws.Range("H122").Value = 4560.9375
$ws.Range("I122").Value = 4497.4
$ws.Range("J122").Value = 4666.8335
$ws.Range("K122").Value = 13492.2
$ws.Range("L122").Value = 14000.5005
$ws.Range("M122").Value = -11042.2
$ws.Range("N122").Value = -18900.5005
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 100000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 100000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 100000
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -100344
$ws.Range("H20").Value = 42510.75
$ws.Range("I20").Value = 10010
$ws.Range("K20").Value = 10010
$ws.Range("M20").Value = -9770
$ws.Range("H54").Value = 10326
$ws.Range("J54").Value = 10326
$ws.Range("L54").Value = 10326
$ws.Range("N54").Value = -11366
$ws.Range("H112").Value = 11346.75
$ws.Range("J112").Value = 11346.75
$ws.Range("L112").Value = 11346.75
$ws.Range("N112").Value = -14300.75
$ws.Range("I132").Value = 2984.7
$ws.Range("J132").Value = 2999.8823
$ws.Range("K132").Value = 8954.099999999999
$ws.Range("L132").Value = 8999.6469
$ws.Range("M132").Value = -6424.099999999999
$ws.Range("N132").Value = -14059.6469
